# Update crypto price (D) and volume(1h) (E) figures with the latest scraped
# data. Price cells that look like plain decimal numbers are entered with a
# leading apostrophe so Excel keeps them as literal text (matching the
# original inline-string formatting, e.g. "1.001", "0.3820") instead of
# silently reinterpreting them as floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.131.15"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.781.83"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'336.56"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.3820"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.3409"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").Value = "'47.91"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").Value = "'1.185"
$ws.Range("E10").Value = "  -4.14%  "
$ws.Range("D11").Value = "'0.07435"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'21.60"
$ws.Range("E13").Value = "  -2.24%  "
$ws.Range("D14").Value = "'6.422"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("D15").Value = "1.783.90"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "'7.085"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").Value = "'0.00001090"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "'0.06635"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").Value = "'83.35"
$ws.Range("E19").Value = "  -3.26%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'6.507"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "'17.34"
$ws.Range("E22").Value = "  -1.61%  "
$ws.Range("D23").Value = "27.134.23"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").Value = "'12.23"
$ws.Range("E24").Value = "  -8.06%  "
$ws.Range("D25").Value = "'2.389"
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("D26").Value = "'2.489"
$ws.Range("E26").Value = "  -7.21%  "
$ws.Range("D27").Value = "'21.07"
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").Value = "'1.440"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").Value = "'154.66"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "1.983.57"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").Value = "'133.64"
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("D32").Value = "'3.980"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").Value = "'6.019"
$ws.Range("E33").Value = "  -5.02%  "
$ws.Range("D34").Value = "'0.08641"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").Value = "'12.94"
$ws.Range("E35").Value = "  -7.10%  "
$ws.Range("D36").Value = "'1.625"
$ws.Range("E36").Value = "  -4.86%  "
$ws.Range("D37").Value = "'5.377"
$ws.Range("E37").Value = "  -4.37%  "
$ws.Range("D38").Value = "'0.6814"
$ws.Range("E38").Value = "  -2.76%  "
$ws.Range("D39").Value = "'0.06285"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").Value = "'0.02325"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("D41").Value = "'0.2171"
$ws.Range("E41").Value = "  -4.78%  "
$ws.Range("D42").Value = "'1.243"
$ws.Range("E42").Value = "  -4.20%  "
$ws.Range("D43").Value = "'8.357"
$ws.Range("E43").Value = "  -7.00%  "
$ws.Range("D44").Value = "'14.17"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'0.6396"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("D47").Value = "'3.855"
$ws.Range("E47").Value = "  -4.88%  "
$ws.Range("D48").Value = "'2.127"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("D49").Value = "'131.18"
$ws.Range("E49").Value = "  -1.42%  "
$ws.Range("D50").Value = "'0.07085"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("D51").Value = "'78.42"
$ws.Range("E51").Value = "  -2.84%  "
